$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Help_text" / "help_desc" columns added to the right of the table (Z, AA)
# describing each field plus an overall table description.

$ws.Range("Z1").Value = "Help_text"
$ws.Range("AA1").Value = "help_desc"

$ws.Range("Z2").Value = "The generic Account code is defined for each entry"

$ws.Range("Z4").Value = "This field denotes the sequence number assigned to an accounting entry in a transaction."

$ws.Range("Z5").Value = "This column represents whether the mentioned Account code has to be debited (+) or credited (-). Both credit and debit amounts for any financial ransaction should tally."

$ws.Range("AA2").Value = "The Accounting rules for all financial transaction codes  need to be configured in this table; these rules  are referred to at the time of processing the respective transactions. Multiple accounting entries could be defined for each transction code. The impacted e Account codes, sequence number and credit(-)/debit(+) details are need to be configured for each entry."

# Highlight the new header cells with a light-blue fill (matches the new fill/style added to the workbook)
$ws.Range("Z1:AA1").Interior.Color = 15773696

# Widen the new help-text column so the long descriptions are readable
$ws.Columns("Z").ColumnWidth = 88.9

# Move the selection to reflect where the edit was made
$ws.Range("Z11").Select()
